$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "release/1.0.1"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"
